# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook

# --- Update timestamps on the "data" sheet (F column) ---
$ws = $wb.Worksheets.Item("data")

$ws.Range("F2").Value = "2021-10-05 14:21:21.089343"
$ws.Range("F3").Value = "2021-10-05 14:21:21.089351"
$ws.Range("F4").Value = "2021-10-05 14:21:21.089354"
$ws.Range("F5").Value = "2021-10-05 14:21:21.089356"
$ws.Range("F6").Value = "2021-10-05 14:21:21.089359"
$ws.Range("F7").Value = "2021-10-05 14:21:21.089361"
$ws.Range("F8").Value = "2021-10-05 14:21:21.089363"
$ws.Range("F9").Value = "2021-10-05 14:21:21.089365"
$ws.Range("F10").Value = "2021-10-05 14:21:21.089367"
$ws.Range("F11").Value = "2021-10-05 14:21:21.089370"
$ws.Range("F12").Value = "2021-10-05 14:21:21.089372"
$ws.Range("F13").Value = "2021-10-05 14:21:21.089374"
$ws.Range("F14").Value = "2021-10-05 14:21:21.089376"

# --- Add a new "metadata" worksheet placed after "data" ---
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Header row (bold, centered, bordered) - clone the exact header style used on
# the "data" sheet (B1:F1) so the new header cells reuse the same style index.
$ws.Range("B1:F1").Copy($meta.Range("B1:F1"))
$ws.Range("B1").Copy($meta.Range("G1"))
# A2 on "data" uses the same bold/bordered style - clone it for metadata!A2 too.
$ws.Range("A2").Copy($meta.Range("A2"))

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row 2
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Kleine-Levin syndrome"
$meta.Range("C2").Value = 213
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.7"
$meta.Range("E2").Value = "2020-10-07T15:03:37.629173Z"
$meta.Range("F2").Value = "2021-10-05 14:21:21.086110"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/213/?format=json"

# Restore "data" as the active sheet/selection (unchanged in the source diff)
$ws.Activate() | Out-Null
$ws.Range("A1").Select() | Out-Null
